$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.788.18"
$ws.Cells.Item(2, 5).Value = "  -0.01%  "
$ws.Cells.Item(3, 4).Value = "1.874.01"
$ws.Cells.Item(3, 5).Value = "  +1.62%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9984"
$ws.Cells.Item(4, 5).Value = "  -0.04%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "242.85"
$ws.Cells.Item(5, 5).Value = "  -2.10%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9981"
$ws.Cells.Item(6, 5).Value = "  -0.08%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4937"
$ws.Cells.Item(7, 5).Value = "  -0.36%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "43.77"
$ws.Cells.Item(8, 5).Value = "  -1.88%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.2897"
$ws.Cells.Item(9, 5).Value = "  +2.95%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.06604"
$ws.Cells.Item(10, 5).Value = "  +2.27%  "
$ws.Cells.Item(11, 4).Value = "1.876.76"
$ws.Cells.Item(11, 5).Value = "  +1.81%  "
$ws.Cells.Item(12, 5).Value = "  +0.09%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.07159"
$ws.Cells.Item(13, 5).Value = "  +0.46%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.6676"
$ws.Cells.Item(14, 5).Value = "  +1.37%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "85.28"
$ws.Cells.Item(15, 5).Value = "  +0.97%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "4.801"
$ws.Cells.Item(16, 5).Value = "  +1.51%  "
$ws.Cells.Item(17, 4).Value = "29.783.19"
$ws.Cells.Item(17, 5).Value = "  +0.08%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000007780"
$ws.Cells.Item(18, 5).Value = "  +5.38%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.9985"
$ws.Cells.Item(19, 5).Value = "  +0.00%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "12.73"
$ws.Cells.Item(20, 5).Value = "  +2.12%  "
$ws.Cells.Item(21, 4).Value = "2.120.72"
$ws.Cells.Item(21, 5).Value = "  +2.25%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.9979"
$ws.Cells.Item(22, 5).Value = "  -0.11%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.726"
$ws.Cells.Item(23, 5).Value = "  +2.79%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "9.126"
$ws.Cells.Item(24, 5).Value = "  +2.67%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "5.560"
$ws.Cells.Item(25, 5).Value = "  +2.39%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "147.23"
$ws.Cells.Item(26, 5).Value = "  +2.72%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "133.89"
$ws.Cells.Item(27, 5).Value = "  +2.04%  "
$ws.Cells.Item(28, 5).Value = "  +1.03%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.918"
$ws.Cells.Item(29, 5).Value = "  +0.94%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.377"
$ws.Cells.Item(30, 5).Value = "  -2.05%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.174"
$ws.Cells.Item(31, 5).Value = "  -0.27%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.08656"
$ws.Cells.Item(32, 5).Value = "  +1.18%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.902"
$ws.Cells.Item(33, 5).Value = "  +2.11%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.05042"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.7060"
$ws.Cells.Item(35, 5).Value = "  +4.59%  "
$ws.Cells.Item(36, 5).Value = "  +0.04%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.667"
$ws.Cells.Item(37, 5).Value = "  -1.15%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.200"
$ws.Cells.Item(38, 5).Value = "  -3.90%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.660"
$ws.Cells.Item(39, 5).Value = "  -1.98%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.9300"
$ws.Cells.Item(40, 5).Value = "  -2.74%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.01634"
$ws.Cells.Item(41, 5).Value = "  +1.95%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "6.054"
$ws.Cells.Item(42, 5).Value = "  -1.14%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.9945"
$ws.Cells.Item(43, 5).Value = "  -0.45%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "102.43"
$ws.Cells.Item(44, 5).Value = "  -0.58%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.4165"
$ws.Cells.Item(45, 5).Value = "  +1.76%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "7.508"
$ws.Cells.Item(46, 5).Value = "  +3.44%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.1254"
$ws.Cells.Item(47, 5).Value = "  +2.10%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.05693"
$ws.Cells.Item(48, 5).Value = "  +1.94%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "32.52"
$ws.Cells.Item(49, 5).Value = "  +1.94%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "8.210"
$ws.Cells.Item(50, 5).Value = "  +1.02%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.337"
$ws.Cells.Item(51, 5).Value = "  +1.47%  "
